$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "date" column (F) values forward by 2 days for rows 2-7
$ws.Range("F2").Value = 44916
$ws.Range("F3").Value = 44915
$ws.Range("F4").Value = 44914
$ws.Range("F5").Value = 44913
$ws.Range("F6").Value = 44912
$ws.Range("F7").Value = 44911
